# Update the cryptocurrency prices/volumes list (GitHub Actions data refresh).
# Rows 50/51 also swap place (dogwifhat <-> FirstDigitalUSD) along with their data.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# stores them as text (matching the original column's text data type).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.533.06"
$ws.Range("E2").Value = "  -2.69%  "
$ws.Range("D3").Value = "3.181.12"
$ws.Range("E3").Value = "  -4.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'570.18"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "'169.18"
$ws.Range("E6").Value = "  -7.86%  "
$ws.Range("E7").Value = "  -5.98%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "3.181.06"
$ws.Range("E9").Value = "  -4.27%  "
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").Value = "'6.81"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("D13").Value = "3.730.68"
$ws.Range("E13").Value = "  -4.43%  "
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "64.571.59"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "'25.41"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "3.176.95"
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("D19").Value = "'420.68"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "'7.12"
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'70.11"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").Value = "'0.206"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'0.489"
$ws.Range("E27").Value = "  -5.34%  "
$ws.Range("E28").Value = "  -7.10%  "
$ws.Range("D29").Value = "'8.88"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -5.76%  "
$ws.Range("D32").Value = "'21.73"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "'5.06"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("E36").Value = "  -4.23%  "
$ws.Range("D37").Value = "'157.77"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("D39").Value = "2.727.26"
$ws.Range("E39").Value = "  -5.52%  "
$ws.Range("E40").Value = "  -5.15%  "
$ws.Range("D41").Value = "'24.35"
$ws.Range("E41").Value = "  -8.12%  "
$ws.Range("D42").Value = "'4.19"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "'39.17"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").Value = "'0.711"
$ws.Range("E44").Value = "  -7.26%  "
$ws.Range("E45").Value = "  -6.34%  "
$ws.Range("D46").Value = "'5.64"
$ws.Range("E46").Value = "  -4.73%  "
$ws.Range("D47").Value = "'0.0264"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").Value = "'294.81"
$ws.Range("E48").Value = "  -6.18%  "
$ws.Range("D49").Value = "'21.65"
$ws.Range("E49").Value = "  -7.09%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'2.01"
$ws.Range("E50").Value = "  -12.97%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.02%  "
